$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-relevant data: rows 2-17, columns A-T

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf2"
$ws.Range("C2").Value = "Sdc3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2151743333333333
$ws.Range("H2").Value = 0.645523
$ws.Range("I2").Value = 0.01945888736810517
$ws.Range("J2").Value = 0.01945888736810517
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 24.44779533333333
$ws.Range("N2").Value = 73.343386
$ws.Range("O2").Value = 0.1890645123346783
$ws.Range("P2").Value = 0.1890645123346783
$ws.Range("Q2").Value = 5.260538062319777
$ws.Range("R2").Value = 47.34484256087799
$ws.Range("S2").Value = 0.003678985050826236
$ws.Range("T2").Value = 0.003678985050826236

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf2"
$ws.Range("C3").Value = "Sdc3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2151743333333333
$ws.Range("H3").Value = 0.645523
$ws.Range("I3").Value = 0.01945888736810517
$ws.Range("J3").Value = 0.01945888736810517
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 18.11074233333333
$ws.Range("N3").Value = 54.332227
$ws.Range("O3").Value = 0.1400575643155068
$ws.Range("P3").Value = 0.1400575643155068
$ws.Range("Q3").Value = 3.896966907746778
$ws.Range("R3").Value = 35.072702169721
$ws.Range("S3").Value = 0.002725364369066593
$ws.Range("T3").Value = 0.002725364369066593

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf2"
$ws.Range("C4").Value = "Sdc3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2151743333333333
$ws.Range("H4").Value = 0.645523
$ws.Range("I4").Value = 0.01945888736810517
$ws.Range("J4").Value = 0.01945888736810517
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 80.684877
$ws.Range("N4").Value = 242.054631
$ws.Range("O4").Value = 0.6239682030546764
$ws.Range("P4").Value = 0.6239682030546765
$ws.Range("Q4").Value = 17.361314618557
$ws.Range("R4").Value = 156.251831567013
$ws.Range("S4").Value = 0.01214172698451992
$ws.Range("T4").Value = 0.01214172698451993

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf2"
$ws.Range("C5").Value = "Sdc3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2151743333333333
$ws.Range("H5").Value = 0.645523
$ws.Range("I5").Value = 0.01945888736810517
$ws.Range("J5").Value = 0.01945888736810517
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.065862
$ws.Range("N5").Value = 18.197586
$ws.Range("O5").Value = 0.0469097202951384
$ws.Range("P5").Value = 0.04690972029513841
$ws.Range("Q5").Value = 1.305217811942
$ws.Range("R5").Value = 11.746960307478
$ws.Range("S5").Value = 0.0009128109636924154
$ws.Range("T5").Value = 0.0009128109636924155

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf2"
$ws.Range("C6").Value = "Sdc3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.077181333333334
$ws.Range("H6").Value = 24.231544
$ws.Range("I6").Value = 0.730444748601188
$ws.Range("J6").Value = 0.730444748601188
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 24.44779533333333
$ws.Range("N6").Value = 73.343386
$ws.Range("O6").Value = 0.1890645123346783
$ws.Range("P6").Value = 0.1890645123346783
$ws.Range("Q6").Value = 197.4692761075538
$ws.Range("R6").Value = 1777.223484967984
$ws.Range("S6").Value = 0.1381011801817103
$ws.Range("T6").Value = 0.1381011801817103

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf2"
$ws.Range("C7").Value = "Sdc3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.077181333333334
$ws.Range("H7").Value = 24.231544
$ws.Range("I7").Value = 0.730444748601188
$ws.Range("J7").Value = 0.730444748601188
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 18.11074233333333
$ws.Range("N7").Value = 54.332227
$ws.Range("O7").Value = 0.1400575643155068
$ws.Range("P7").Value = 0.1400575643155068
$ws.Range("Q7").Value = 146.2837499076098
$ws.Range("R7").Value = 1316.553749168488
$ws.Range("S7").Value = 0.1023043123561351
$ws.Range("T7").Value = 0.1023043123561351

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fgf2"
$ws.Range("C8").Value = "Sdc3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.077181333333334
$ws.Range("H8").Value = 24.231544
$ws.Range("I8").Value = 0.730444748601188
$ws.Range("J8").Value = 0.730444748601188
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 80.684877
$ws.Range("N8").Value = 242.054631
$ws.Range("O8").Value = 0.6239682030546764
$ws.Range("P8").Value = 0.6239682030546765
$ws.Range("Q8").Value = 651.7063823866961
$ws.Range("R8").Value = 5865.357441480264
$ws.Range("S8").Value = 0.4557742972154081
$ws.Range("T8").Value = 0.4557742972154082

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fgf2"
$ws.Range("C9").Value = "Sdc3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.077181333333334
$ws.Range("H9").Value = 24.231544
$ws.Range("I9").Value = 0.730444748601188
$ws.Range("J9").Value = 0.730444748601188
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.065862
$ws.Range("N9").Value = 18.197586
$ws.Range("O9").Value = 0.0469097202951384
$ws.Range("P9").Value = 0.04690972029513841
$ws.Range("Q9").Value = 48.99506731697601
$ws.Range("R9").Value = 440.955605852784
$ws.Range("S9").Value = 0.03426495884793442
$ws.Range("T9").Value = 0.03426495884793442

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Fgf2"
$ws.Range("C10").Value = "Sdc3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.519176
$ws.Range("H10").Value = 1.557528
$ws.Range("I10").Value = 0.04695070806875992
$ws.Range("J10").Value = 0.04695070806875992
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 24.44779533333333
$ws.Range("N10").Value = 73.343386
$ws.Range("O10").Value = 0.1890645123346783
$ws.Range("P10").Value = 0.1890645123346783
$ws.Range("Q10").Value = 12.69270858997866
$ws.Range("R10").Value = 114.234377309808
$ws.Range("S10").Value = 0.00887671272478794
$ws.Range("T10").Value = 0.008876712724787942

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Fgf2"
$ws.Range("C11").Value = "Sdc3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.519176
$ws.Range("H11").Value = 1.557528
$ws.Range("I11").Value = 0.04695070806875992
$ws.Range("J11").Value = 0.04695070806875992
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 18.11074233333333
$ws.Range("N11").Value = 54.332227
$ws.Range("O11").Value = 0.1400575643155068
$ws.Range("P11").Value = 0.1400575643155068
$ws.Range("Q11").Value = 9.402662761650667
$ws.Range("R11").Value = 84.62396485485601
$ws.Range("S11").Value = 0.006575801814998928
$ws.Range("T11").Value = 0.006575801814998929

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Fgf2"
$ws.Range("C12").Value = "Sdc3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.519176
$ws.Range("H12").Value = 1.557528
$ws.Range("I12").Value = 0.04695070806875992
$ws.Range("J12").Value = 0.04695070806875992
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 80.684877
$ws.Range("N12").Value = 242.054631
$ws.Range("O12").Value = 0.6239682030546764
$ws.Range("P12").Value = 0.6239682030546765
$ws.Range("Q12").Value = 41.889651701352
$ws.Range("R12").Value = 377.006865312168
$ws.Range("S12").Value = 0.02929574894580882
$ws.Range("T12").Value = 0.02929574894580883

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Fgf2"
$ws.Range("C13").Value = "Sdc3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.519176
$ws.Range("H13").Value = 1.557528
$ws.Range("I13").Value = 0.04695070806875992
$ws.Range("J13").Value = 0.04695070806875992
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 6.065862
$ws.Range("N13").Value = 18.197586
$ws.Range("O13").Value = 0.0469097202951384
$ws.Range("P13").Value = 0.04690972029513841
$ws.Range("Q13").Value = 3.149249969712
$ws.Range("R13").Value = 28.343249727408
$ws.Range("S13").Value = 0.002202444583164225
$ws.Range("T13").Value = 0.002202444583164226

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Fgf2"
$ws.Range("C14").Value = "Sdc3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.246363333333333
$ws.Range("H14").Value = 6.73909
$ws.Range("I14").Value = 0.2031456559619469
$ws.Range("J14").Value = 0.2031456559619469
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 24.44779533333333
$ws.Range("N14").Value = 73.343386
$ws.Range("O14").Value = 0.1890645123346783
$ws.Range("P14").Value = 0.1890645123346783
$ws.Range("Q14").Value = 54.91863101763778
$ws.Range("R14").Value = 494.26767915874
$ws.Range("S14").Value = 0.03840763437735383
$ws.Range("T14").Value = 0.03840763437735384

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Fgf2"
$ws.Range("C15").Value = "Sdc3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.246363333333333
$ws.Range("H15").Value = 6.73909
$ws.Range("I15").Value = 0.2031456559619469
$ws.Range("J15").Value = 0.2031456559619469
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 18.11074233333333
$ws.Range("N15").Value = 54.332227
$ws.Range("O15").Value = 0.1400575643155068
$ws.Range("P15").Value = 0.1400575643155068
$ws.Range("Q15").Value = 40.68330751704778
$ws.Range("R15").Value = 366.14976765343
$ws.Range("S15").Value = 0.02845208577530621
$ws.Range("T15").Value = 0.02845208577530621

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Fgf2"
$ws.Range("C16").Value = "Sdc3"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.246363333333333
$ws.Range("H16").Value = 6.73909
$ws.Range("I16").Value = 0.2031456559619469
$ws.Range("J16").Value = 0.2031456559619469
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 80.684877
$ws.Range("N16").Value = 242.054631
$ws.Range("O16").Value = 0.6239682030546764
$ws.Range("P16").Value = 0.6239682030546765
$ws.Range("Q16").Value = 181.24754924731
$ws.Range("R16").Value = 1631.22794322579
$ws.Range("S16").Value = 0.1267564299089395
$ws.Range("T16").Value = 0.1267564299089396

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Fgf2"
$ws.Range("C17").Value = "Sdc3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.246363333333333
$ws.Range("H17").Value = 6.73909
$ws.Range("I17").Value = 0.2031456559619469
$ws.Range("J17").Value = 0.2031456559619469
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 6.065862
$ws.Range("N17").Value = 18.197586
$ws.Range("O17").Value = 0.0469097202951384
$ws.Range("P17").Value = 0.04690972029513841
$ws.Range("Q17").Value = 13.62612998186
$ws.Range("R17").Value = 122.63516983674
$ws.Range("S17").Value = 0.009529505900347346
$ws.Range("T17").Value = 0.009529505900347348
